$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("optimal models")
$ws2 = $wb.Worksheets.Item("relative importance")

# --- Sheet1: "optimal models" - add row 8 for BET ---
$ws1.Range("A8").Value = "BET"
$ws1.Range("B8").Value = 2
$ws1.Range("C8").Value = 0.75
$ws1.Range("D8").Value = 0.006
$ws1.Range("E8").Value = 6600
$ws1.Range("F8").Value = 0.9462
$ws1.Range("G8").Value = 0.86532
$ws1.Range("H8").Value = 0.8194
$ws1.Range("I8").Value = 0.2027264

# --- Sheet2: "relative importance" - add row 8 for BET ---
$ws2.Range("A8").Value = "BET"
$ws2.Range("B8").Value = 9.959357
$ws2.Range("C8").Value = 3.206406
$ws2.Range("D8").Value = 4.529701
$ws2.Range("E8").Value = 3.162866
$ws2.Range("F8").Value = 6.142197
$ws2.Range("G8").Value = 10.81284
$ws2.Range("H8").Value = 2.179942
$ws2.Range("I8").Value = 4.62519
$ws2.Range("J8").Value = 5.369122
$ws2.Range("K8").Value = 2.399092
$ws2.Range("L8").Value = 7.799043
$ws2.Range("M8").Value = 1.399017
$ws2.Range("N8").Value = 11.430529
$ws2.Range("O8").Value = 10.927342
$ws2.Range("P8").Value = 6.366438
$ws2.Range("Q8").Value = 6.462402
$ws2.Range("R8").Value = 3.228516

# --- Update selections to match target state ---
# Select sheet2's cell first, then sheet1's cell last so sheet1 ends up
# as the active tab (matches tabSelected="1" staying on sheet1).
$ws2.Range("M9").Select()
$ws1.Range("H3").Select()
